$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.964.86"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.60%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.92"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.18%  "

# Row 4
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.21%  "

# Row 6
$ws.Range("E6").Value = "  -0.29%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4568"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.33%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3698"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.84%  "

# Row 9
$ws.Range("E9").Value = "  +0.06%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8719"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.09%  "

# Row 11
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.963.64"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +7.06%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07958"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.41%  "

# Row 13
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.73"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.97%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.547"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.09%  "

# Row 15
$ws.Range("E15").Value = "  -0.43%  "

# Row 16
$ws.Range("E16").Value = "  -1.28%  "

# Row 17
$ws.Range("E17").Value = "  -0.18%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008853"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.29%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.22%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.70"
$ws.Range("D20").ClearFormats()

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.121.03"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.13%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.099"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.26%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.52"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.33%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.148.18"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.26"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.846"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.26%  "

# Row 27
$ws.Range("E27").Value = "  +0.99%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.037"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.98%  "

# Row 29
$ws.Range("E29").Value = "  +0.87%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.04"
$ws.Range("D30").ClearFormats()

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08873"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.43%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.965"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.19%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7270"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.25%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.410"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.04%  "

# Row 35
$ws.Range("E35").Value = "  -1.05%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.074"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.30%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.441"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.06%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01937"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05219"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.939"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.52%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.159"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.09%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5128"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.54%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1625"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.25%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.168"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.46%  "

# Row 45
$ws.Range("E45").Value = "  -0.42%  "

# Row 46
$ws.Range("E46").Value = "  -0.27%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.15"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.03%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.35"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.00%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.628"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.50%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06206"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.88%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.04"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.86%  "
